# Update countries & provincias Spain
# Applies the data refresh described by the commit "Update countries & provincias Spain"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Updated timestamp banner (A1)
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 22 de Julio de 2020 a las 23:15"

# Israel overtakes Republica Dominicana in total cases -> rows 41/42 swap countries,
# with Israel's stats refreshed to the new totals and Republica Dominicana's stats
# staying the same (just moved down one row).
$ws.Cells.Item(41, 1).Value = "Israel"
$ws.Cells.Item(41, 2).Value = 56085
$ws.Cells.Item(41, 3).Value = 2043
$ws.Cells.Item(41, 4).Value = 23310
$ws.Cells.Item(41, 5).Value = 32345
$ws.Cells.Item(41, 6).Value = 0
$ws.Cells.Item(41, 7).Value = 5
$ws.Cells.Item(41, 8).Value = 430

$ws.Cells.Item(42, 1).Value = "Republica Dominicana"
$ws.Cells.Item(42, 2).Value = 56043
$ws.Cells.Item(42, 3).Value = 1246
$ws.Cells.Item(42, 4).Value = 26466
$ws.Cells.Item(42, 5).Value = 28572
$ws.Cells.Item(42, 6).Value = 0
$ws.Cells.Item(42, 7).Value = 6
$ws.Cells.Item(42, 8).Value = 1005

# Groenlandia and Islas Malvinas swap positions (tied totals, labels only)
$ws.Cells.Item(210, 1).Value = "Groenlandia"
$ws.Cells.Item(211, 1).Value = "Islas Malvinas"

# Refreshed case counts for several countries
$ws.Cells.Item(4, 2).Value = 4076721
$ws.Cells.Item(4, 3).Value = 48152
$ws.Cells.Item(4, 4).Value = 1916197
$ws.Cells.Item(4, 5).Value = 2014738
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(4, 7).Value = 833
$ws.Cells.Item(4, 8).Value = 145786

$ws.Cells.Item(5, 2).Value = 2227514
$ws.Cells.Item(5, 3).Value = 60982
$ws.Cells.Item(5, 4).Value = 1465970
$ws.Cells.Item(5, 5).Value = 678773
$ws.Cells.Item(5, 6).Value = 0
$ws.Cells.Item(5, 7).Value = 1174
$ws.Cells.Item(5, 8).Value = 82771

$ws.Cells.Item(8, 2).Value = 394948
$ws.Cells.Item(8, 3).Value = 13150
$ws.Cells.Item(8, 4).Value = 229175
$ws.Cells.Item(8, 5).Value = 159833
$ws.Cells.Item(8, 6).Value = 0
$ws.Cells.Item(8, 7).Value = 572
$ws.Cells.Item(8, 8).Value = 5940

$ws.Cells.Item(51, 2).Value = 37637
$ws.Cells.Item(51, 3).Value = 321
$ws.Cells.Item(51, 4).Value = 33894
$ws.Cells.Item(51, 5).Value = 3613

$ws.Cells.Item(65, 2).Value = 18379
$ws.Cells.Item(65, 3).Value = 498
$ws.Cells.Item(65, 4).Value = 9872
$ws.Cells.Item(65, 5).Value = 8409
$ws.Cells.Item(65, 6).Value = 0
$ws.Cells.Item(65, 7).Value = 3
$ws.Cells.Item(65, 8).Value = 98
